$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores every Coin/Link/Price/Volume cell as plain
# text (inline strings), including price strings that look numeric
# (e.g. "298.01") and multi-dot big-number strings (e.g. "42.140.82").
# Force text number-format before writing so Excel does not silently
# re-interpret these as numbers, then clear the temporary format so
# the cell style is left exactly as it was (no lingering "@" style).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.140.82"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.98%  "
$ws.Range("E2").ClearFormats()
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.269.21"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.88%  "
$ws.Range("E3").ClearFormats()
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E4").ClearFormats()
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.01"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.78%  "
$ws.Range("E5").ClearFormats()
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.31"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.61%  "
$ws.Range("E6").ClearFormats()
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.494"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.64%  "
$ws.Range("E7").ClearFormats()
# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E8").ClearFormats()
# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.95%  "
$ws.Range("E9").ClearFormats()
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.02"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.12%  "
$ws.Range("E10").ClearFormats()
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("E11").ClearFormats()
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.10"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -7.99%  "
$ws.Range("E12").ClearFormats()
# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("E13").ClearFormats()
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.65"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.82%  "
$ws.Range("E14").ClearFormats()
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.61"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("E15").ClearFormats()
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.620.90"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.06%  "
$ws.Range("E16").ClearFormats()
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.254.94"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.66%  "
$ws.Range("E17").ClearFormats()
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.772"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.04%  "
$ws.Range("E18").ClearFormats()
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.130.68"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("E19").ClearFormats()
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.49%  "
$ws.Range("E20").ClearFormats()
# Row 21
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("B21").ClearFormats()
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C21").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.33"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.55%  "
$ws.Range("E21").ClearFormats()
# Row 22
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "Uniswap"
$ws.Range("B22").ClearFormats()
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("C22").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.97"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.07%  "
$ws.Range("E22").ClearFormats()
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.68"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.75%  "
$ws.Range("E23").ClearFormats()
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "233.06"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.67%  "
$ws.Range("E24").ClearFormats()
# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.11%  "
$ws.Range("E25").ClearFormats()
# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E26").ClearFormats()
# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.18%  "
$ws.Range("E27").ClearFormats()
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.81"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.33%  "
$ws.Range("E28").ClearFormats()
# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("E29").ClearFormats()
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.57"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("E30").ClearFormats()
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.62"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.83%  "
$ws.Range("E31").ClearFormats()
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.03"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.80%  "
$ws.Range("E32").ClearFormats()
# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E33").ClearFormats()
# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.08%  "
$ws.Range("E34").ClearFormats()
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.47"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("E35").ClearFormats()
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.33"
$ws.Range("D36").ClearFormats()
# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.14%  "
$ws.Range("E37").ClearFormats()
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.03"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -8.83%  "
$ws.Range("E38").ClearFormats()
# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.02%  "
$ws.Range("E39").ClearFormats()
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0988"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.68%  "
$ws.Range("E40").ClearFormats()
# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("E41").ClearFormats()
# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.94%  "
$ws.Range("E42").ClearFormats()
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.45"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("E43").ClearFormats()
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.961.42"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.05%  "
$ws.Range("E44").ClearFormats()
# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.03%  "
$ws.Range("E45").ClearFormats()
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.25"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.73%  "
$ws.Range("E46").ClearFormats()
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.53"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.23%  "
$ws.Range("E47").ClearFormats()
# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.26%  "
$ws.Range("E48").ClearFormats()
# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("B49").ClearFormats()
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("C49").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.493.29"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.49%  "
$ws.Range("E49").ClearFormats()
# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("B50").ClearFormats()
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C50").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.79"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.39%  "
$ws.Range("E50").ClearFormats()
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.96"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.99%  "
$ws.Range("E51").ClearFormats()
